$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.722348213195801
$ws.Range("B1").Value = 2.909647226333618
$ws.Range("C1").Value = 1.986009120941162
$ws.Range("D1").Value = 1.623988509178162
$ws.Range("E1").Value = 1.51916778087616
